$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values must stay as text (they mimic numbers like "0.9999" or
# contain multiple dots like "28.024.46"), so force text format before
# assigning, then restore the default "Normal" style so no stray cell
# style is left behind.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.024.46'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.10%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.872.65'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.73%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.33%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.61'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.03%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9994'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.29%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5043'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3828'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.45%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09013'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -5.80%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.117'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.26%  '

$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.67'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.96%  '

$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.371'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.07%  '

$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.70'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.08%  '

$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.869.79'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.56%  '

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.256'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.42%  '

$ws.Range('B16').Value = 'BinanceUSD'
$ws.Range('C16').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9987'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.45%  '

$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001106'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.79%  '

$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '91.16'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.05%  '

$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06643'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.77%  '

$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.15'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.73%  '

$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9990'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.27%  '

$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.121'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.35%  '

$ws.Range('B23').Value = 'WrappedBTC'
$ws.Range('C23').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.048.59'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.29%  '

$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.51'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.41%  '

$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.256'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.07%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.079.26'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.48%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.518'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.09%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.75'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.13%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '156.99'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.46%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.89'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.29%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1066'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.51%  '

$ws.Range('E32').Value = '  -2.32%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.602'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.10%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.592'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.71%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.453'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.54%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06577'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.44%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02402'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.71%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2192'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.67%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.291'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.10%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.207'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.00%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6388'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.93%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.47'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.89%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.913'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.18%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9986'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.33%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.24'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.55%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6021'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.67%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.274'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.37%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.661'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.81%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.246'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.61%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.996'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.17%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.78'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.54%  '
